$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix typo: "generell" -> "general" for biodiv_invest_objective (row 57)
$ws.Range("F57").Value = '"general"'

# 2. Replace the old "impact_all" parameter row (row 58) with the new
#    "which-machine?" parameter (string, default "server")
$ws.Range("A58").Value = "which-machine?"
$ws.Range("E58").Value = "string"
$ws.Range("F58").Value = '"server"'

# 3. Remove now-obsolete variables. Delete from the bottom up so earlier
#    row numbers stay valid while we work.
$ws.Rows("111:111").Delete()   # p_habitat_quality
$ws.Rows("88:89").Delete()     # landscape-hq, dist_max
